$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right after "2021-Q3" and before "总计"
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q3")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "2022-Q1"

# Copy the header/index-cell formatting (style "s=2": bold, centered, boxed)
# from the existing "2021-Q3" sheet so the new sheet matches the established
# look used by the other quarterly sheets.
$wb.Worksheets.Item("2021-Q3").Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$wb.Worksheets.Item("2021-Q3").Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row. Columns B-G hold text (even the numeric-looking ones), matching
# the convention used throughout the other quarterly sheets; H holds a real
# number. Flip each text cell's number format to Text before assigning the
# value so it's stored as a string, then strip the format again afterwards
# so no stray style index is left on the cell.
$newSheet.Range("A2").Value = 0

$textCells = "B2", "C2", "D2", "E2", "F2", "G2"
foreach ($addr in $textCells) {
    $newSheet.Range($addr).NumberFormat = "@"
}
$newSheet.Range("B2").Value = "160125"
$newSheet.Range("C2").Value = "南方香港优选股票QDII-LOF"
$newSheet.Range("D2").Value = "2.46"
$newSheet.Range("E2").Value = "91.14"
$newSheet.Range("F2").Value = "3.67"
$newSheet.Range("G2").Value = "0.0903"
foreach ($addr in $textCells) {
    $newSheet.Range($addr).ClearFormats()
}
$newSheet.Range("H2").Value = 4

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new first data row for
#    2022-Q1, shifting the existing quarterly summary rows down by one and
#    renumbering their index column (A).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.09

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

# Restore the originally active sheet/selection (adding a sheet shifts focus
# to it by default).
$wb.Worksheets.Item("2020-Q4").Activate()
[void]$wb.Worksheets.Item("2020-Q4").Range("A1").Select()
